$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.874.58'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +4.91%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.277.68'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +5.27%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '580.84'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +2.56%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '182.43'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +8.61%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.601'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.275.99'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +5.24%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.134'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +9.47%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.75'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +3.52%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.419'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +8.27%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.849.08'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +5.33%  '
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +1.29%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '28.62'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +7.76%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '67.852.61'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +5.02%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000169'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +5.43%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.279.50'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +5.12%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.86'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +3.84%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.59'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +7.73%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '375.75'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +6.60%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.67'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +7.41%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.20%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '71.23'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +3.78%  '
$ws.Range("B25").Value = 'Polygon'
$ws.Range("C25").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.513'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +4.88%  '
$ws.Range("B26").Value = 'PEPE'
$ws.Range("C26").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000120'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +6.66%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.66'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +1.21%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +3.74%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.15%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.98'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +4.56%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.70'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +9.78%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '22.79'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +5.54%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.28'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +8.95%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.04%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +6.74%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.51'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +7.10%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '163.51'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +3.39%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.852'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +3.61%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +7.08%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +13.37%  '
$ws.Range("B41").Value = 'Filecoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.68'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +13.57%  '
$ws.Range("B42").Value = 'EnergySwap'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '26.86'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +3.39%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.63'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +10.38%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '355.82'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +12.68%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.712.23'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +3.33%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '25.57'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +8.43%  '
$ws.Range("B47").Value = 'OKB'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '40.90'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +4.22%  '
$ws.Range("B48").Value = 'Hedera'
$ws.Range("C48").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0683'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +5.48%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +4.82%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +8.32%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.99%  '
